$d = $word.ActiveDocument
$d.Content.Find.Execute("50ml sterile bottle to use", $true, $false, $false, $false, $false,
                         $true, 1, $false, "50ml sterile falcon tube to use", 2)
